$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header values (B1:E1) ---
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# --- Row 2 values ---
$ws.Range("B2").Value = 1.7052526527940211
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 2.2283496212090781
$ws.Range("E2").ClearContents()

# --- Row 3 values ---
$ws.Range("B3").Value = 1.540892365117235
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 1.9469526140131026
$ws.Range("E3").Value = -1.7914419616663402

# --- Selection range updated to reflect new used range ---
$ws.Range("B1:E3").Select() | Out-Null
